$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModuleController")

# Add a new data row (row 9) to the ModuleController list:
#   ExecutionFlag (col A) = "Yes", ModuleName (col B) = "FRAT"
$ws.Range("A9").Value = "Yes"
$ws.Range("B9").Value = "FRAT"

# Update the active selection, as recorded when the workbook was last saved
$ws.Range("A13").Select()
